$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("20151105")

$ws.Range("A3").Value = "Add support for screen size of 632 by 1030"
$ws.Range("B3").Value = "Feature Request"

$ws.Range("D3").Select()
